$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 921.05
$ws.Range("I28").Value = 578.625
$ws.Range("J28").Value = 2290.75
$ws.Range("K28").Value = 578.625
$ws.Range("L28").Value = 2290.75
$ws.Range("M28").Value = -93.625
$ws.Range("N28").Value = -3260.75
$ws.Range("H138").Value = 4755.4165
$ws.Range("I138").Value = 4824.091
$ws.Range("J138").Value = 4000
$ws.Range("K138").Value = 14472.273
$ws.Range("L138").Value = 12000
$ws.Range("M138").Value = -9332.273000000001
$ws.Range("N138").Value = -22280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4764357
$ws.Range("I61").Value = 2499.6843
$ws.Range("K61").Value = 2499.6843
$ws.Range("M61").Value = -2287.6843
$ws.Range("H97").Value = 1540.5
$ws.Range("I97").Value = 1540.5
$ws.Range("K97").Value = 1540.5
$ws.Range("M97").Value = -1044.5
$ws.Range("H124").Value = 18678
$ws.Range("J124").Value = 18678
$ws.Range("L124").Value = 18678
$ws.Range("N124").Value = -28498
$ws.Range("H131").Value = 50000
$ws.Range("J131").Value = 50000
$ws.Range("L131").Value = 50000
$ws.Range("N131").Value = -60080
$ws.Range("H136").Value = 4764357
$ws.Range("I136").Value = 2499.6843
$ws.Range("K136").Value = 7499.0529
$ws.Range("M136").Value = -4949.0529

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 2050.75
$ws.Range("J10").Value = 2940.6
$ws.Range("L10").Value = 2940.6
$ws.Range("N10").Value = -3220.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 616.9666999999999
$ws.Range("I10").Value = 626.2593000000001
$ws.Range("J10").Value = 533.3333
$ws.Range("K10").Value = 626.2593000000001
$ws.Range("L10").Value = 533.3333
$ws.Range("M10").Value = -487.2593000000001
$ws.Range("N10").Value = -811.3333
$ws.Range("H20").Value = 60970
$ws.Range("I20").Value = 10000
$ws.Range("J20").Value = 73712.5
$ws.Range("K20").Value = 10000
$ws.Range("L20").Value = 73712.5
$ws.Range("M20").Value = -9764
$ws.Range("N20").Value = -74184.5
$ws.Range("H22").Value = 1387.8695
$ws.Range("I22").Value = 977.61536
$ws.Range("J22").Value = 1921.2
$ws.Range("K22").Value = 977.61536
$ws.Range("L22").Value = 1921.2
$ws.Range("M22").Value = -627.61536
$ws.Range("N22").Value = -2621.2
$ws.Range("H30").Value = 60970
$ws.Range("I30").Value = 10000
$ws.Range("J30").Value = 73712.5
$ws.Range("K30").Value = 10000
$ws.Range("L30").Value = 73712.5
$ws.Range("M30").Value = -9909
$ws.Range("N30").Value = -73894.5
$ws.Range("H31").Value = 903926.9
$ws.Range("I31").Value = 1230518.5
$ws.Range("J31").Value = 5799.75
$ws.Range("K31").Value = 1230518.5
$ws.Range("L31").Value = 5799.75
$ws.Range("M31").Value = -1230223.5
$ws.Range("N31").Value = -6389.75
$ws.Range("H34").Value = 903926.9
$ws.Range("I34").Value = 1230518.5
$ws.Range("J34").Value = 5799.75
$ws.Range("K34").Value = 1230518.5
$ws.Range("L34").Value = 5799.75
$ws.Range("M34").Value = -1230316.5
$ws.Range("N34").Value = -6203.75
$ws.Range("H122").Value = 10880
$ws.Range("I122").Value = 1578.6
$ws.Range("J122").Value = 21612.385
$ws.Range("K122").Value = 4735.799999999999
$ws.Range("L122").Value = 64837.155
$ws.Range("M122").Value = -2285.799999999999
$ws.Range("N122").Value = -69737.155
$ws.Range("H128").Value = 60970
$ws.Range("I128").Value = 10000
$ws.Range("J128").Value = 73712.5
$ws.Range("K128").Value = 10000
$ws.Range("L128").Value = 73712.5
$ws.Range("M128").Value = -5020
$ws.Range("N128").Value = -83672.5
$ws.Range("H131").Value = 59849.5
$ws.Range("J131").Value = 59849.5
$ws.Range("L131").Value = 59849.5
$ws.Range("N131").Value = -69929.5
$ws.Range("H132").Value = 2821.8948
$ws.Range("I132").Value = 2483.2144
$ws.Range("K132").Value = 7449.6432
$ws.Range("M132").Value = -4919.6432
$ws.Range("H134").Value = 3974.75
$ws.Range("I134").Value = 3135.9092
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 9407.7276
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -6872.7276
$ws.Range("N134").Value = -20070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 90910920
$ws.Range("I11").Value = 5000.75
$ws.Range("K11").Value = 15002.25
$ws.Range("M11").Value = -14862.25
$ws.Range("H68").Value = 5092
$ws.Range("I68").Value = 1916.6666
$ws.Range("J68").Value = 6044.6
$ws.Range("K68").Value = 5749.9998
$ws.Range("L68").Value = 18133.8
$ws.Range("M68").Value = -4938.9998
$ws.Range("N68").Value = -19755.8
$ws.Range("H70").Value = 7254.8
$ws.Range("I70").Value = 2166.3333
$ws.Range("K70").Value = 6498.999899999999
$ws.Range("M70").Value = -6183.999899999999
$ws.Range("H71").Value = 5092
$ws.Range("I71").Value = 1916.6666
$ws.Range("J71").Value = 6044.6
$ws.Range("K71").Value = 17249.9994
$ws.Range("L71").Value = 54401.4
$ws.Range("M71").Value = -13193.9994
$ws.Range("N71").Value = -62513.4
$ws.Range("H73").Value = 7254.8
$ws.Range("I73").Value = 2166.3333
$ws.Range("K73").Value = 6498.999899999999
$ws.Range("M73").Value = -5406.999899999999
$ws.Range("H108").Value = 1950
$ws.Range("I108").Value = 933.3333
$ws.Range("K108").Value = 2799.9999
$ws.Range("M108").Value = 80.0001000000002
$ws.Range("H109").Value = 6123.2856
$ws.Range("I109").Value = 431.75
$ws.Range("K109").Value = 1295.25
$ws.Range("M109").Value = -255.25
$ws.Range("H119").Value = 9174.532999999999
$ws.Range("I119").Value = 3327.25
$ws.Range("J119").Value = 15857.143
$ws.Range("K119").Value = 9981.75
$ws.Range("L119").Value = 47571.429
$ws.Range("M119").Value = -5143.75
$ws.Range("N119").Value = -57247.429
$ws.Range("H121").Value = 3098.6191
$ws.Range("J121").Value = 3418.9443
$ws.Range("L121").Value = 10256.8329
$ws.Range("N121").Value = -12876.8329
$ws.Range("H138").Value = 18996.344
$ws.Range("I138").Value = 27322.666
$ws.Range("K138").Value = 81967.99800000001
$ws.Range("M138").Value = -76827.99800000001
$ws.Range("H140").Value = 2605
$ws.Range("I140").Value = 2238.75
$ws.Range("J140").Value = 7000
$ws.Range("K140").Value = 6716.25
$ws.Range("L140").Value = 21000
$ws.Range("M140").Value = -1536.25
$ws.Range("N140").Value = -31360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 5010751
$ws.Range("I18").Value = 5010751
$ws.Range("K18").Value = 5010751
$ws.Range("M18").Value = -5010458
$ws.Range("H21").Value = 129356.6
$ws.Range("I21").Value = 3142.7144
$ws.Range("J21").Value = 423855.66
$ws.Range("K21").Value = 3142.7144
$ws.Range("L21").Value = 423855.66
$ws.Range("M21").Value = -2969.7144
$ws.Range("N21").Value = -424201.66
$ws.Range("H30").Value = 129356.6
$ws.Range("I30").Value = 3142.7144
$ws.Range("J30").Value = 423855.66
$ws.Range("K30").Value = 3142.7144
$ws.Range("L30").Value = 423855.66
$ws.Range("M30").Value = -3037.7144
$ws.Range("N30").Value = -424065.66
$ws.Range("H43").Value = 18500
$ws.Range("I43").Value = 12000
$ws.Range("K43").Value = 12000
$ws.Range("M43").Value = -11849
$ws.Range("H129").Value = 56500
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 56500
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 56500
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -66500

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 3671.5
$ws.Range("I9").Value = 507.25
$ws.Range("K9").Value = 507.25
$ws.Range("M9").Value = -283.25
$ws.Range("H23").Value = 11994.4
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H30").Value = 830.4
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()
$ws.Range("H130").Value = 49333
$ws.Range("J130").Value = 49333
$ws.Range("L130").Value = 49333
$ws.Range("N130").Value = -59373
$ws.Range("H131").Value = 57764.5
$ws.Range("J131").Value = 57764.5
$ws.Range("L131").Value = 57764.5
$ws.Range("N131").Value = -67844.5
